$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 29 de Junio de 2020 a las 08:01"

# Swap Fiyi/Dominica order: row 205 becomes Dominica, row 206 becomes Fiyi
# (their numeric data is identical, so swapping the labels matches the
# shared-string reordering in the diff)
$ws.Range("A205").Value = "Dominica"
$ws.Range("A206").Value = "Fiyi"

# Row 47: Afganistan
$ws.Range("B47").Value = 31238
$ws.Range("C47").Value = 271
$ws.Range("D47").Value = 13934
$ws.Range("E47").Value = 16571
$ws.Range("G47").Value = 12
$ws.Range("H47").Value = 733

# Row 74: Uzbekistan
$ws.Range("B74").Value = 8031
$ws.Range("C74").Value = 83
$ws.Range("E74").Value = 2680

# Row 75: Australia
$ws.Range("B75").Value = 7767
$ws.Range("C75").Value = 81
$ws.Range("D75").Value = 7008
$ws.Range("E75").Value = 655

# Row 81: El Salvador
$ws.Range("D81").Value = 3566
$ws.Range("E81").Value = 2204
$ws.Range("G81").Value = 12
$ws.Range("H81").Value = 164

# Row 98: Tailandia
$ws.Range("B98").Value = 3169
$ws.Range("C98").Value = 7
$ws.Range("E98").Value = 58
